$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns at I:J. This pushes the old I,J,K,L,M,N content
#    (Created By / Date Created labels + inputs, and the trailing filler
#    columns) two slots to the right -> K,L,M,N,O,P. It also naturally grows
#    the B2:J2 merged banner to B2:L2 and re-targets the border styling used
#    on that merged row, matching the target workbook exactly.
# ---------------------------------------------------------------------------
$ws.Columns("I:J").Insert()

# ---------------------------------------------------------------------------
# 2. Re-arrange the results-table header (row 8). "Result" moves after the
#    Fine/Gross/Others columns, and two brand-new columns "Defect No" /
#    "Defect Desc." are introduced right before "Created By".
# ---------------------------------------------------------------------------
$ws.Range("E8").Value = "Fine"
$ws.Range("F8").Value = "Gross"
$ws.Range("G8").Value = "Others"
$ws.Range("H8").Value = "Result"
$ws.Range("I8").Value = "Defect No"
$ws.Range("J8").Value = "Defect Desc."

# ---------------------------------------------------------------------------
# 3. Column G in the "Generated Date / Generated By" rows (5 & 6) picks up
#    the same left/vertical-centred alignment the neighbouring filler cells
#    already use.
# ---------------------------------------------------------------------------
$ws.Range("G5:G6").HorizontalAlignment = -4131
$ws.Range("G5:G6").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Resize the (now six) data columns E:J that hold the per-record detail
#    fields to their new, final widths.
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 11.83
$ws.Columns("F").ColumnWidth = 8.5
$ws.Columns("G").ColumnWidth = 11.67
$ws.Columns("H").ColumnWidth = 7.33
$ws.Columns("I").ColumnWidth = 12.67
$ws.Columns("J").ColumnWidth = 21

# ---------------------------------------------------------------------------
# 5. Restore the last-used selection, which Excel persisted as L5 after the
#    edit (the "Generated Date" input cell, post column-insert).
# ---------------------------------------------------------------------------
$ws.Range("L5").Select()

# ---------------------------------------------------------------------------
# 6. Sharepoint sync path recorded by the desktop client that last saved the
#    file moved from "David's Stuff" to "Desktop".
# ---------------------------------------------------------------------------
try {
    $wb.AbsPath = "https://johnsonoutdoorsinc-my.sharepoint.com/personal/david_zhang_johnsonoutdoors_com/Documents/Desktop/"
} catch {
}
